$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The whole "2020" column (B) is dropped and every other year column shifts
# one slot to the left (2021->B, 2022->C, 2023->D, 2024->E, 2025->F), while
# the underlying monthly totals were refreshed from the latest source data.
$ws.Columns("B").Delete()

# Header row (years)
$ws.Range("B1").Value = 2021
$ws.Range("C1").Value = 2022
$ws.Range("D1").Value = 2023
$ws.Range("E1").Value = 2024
$ws.Range("F1").Value = 2025

# Refreshed monthly totals
$ws.Range("B2").Value = 50518.99
$ws.Range("C2").Value = 141789.57
$ws.Range("D2").Value = 162810.55
$ws.Range("E2").Value = 331129.46
$ws.Range("F2").Value = 402689.42

$ws.Range("B3").Value = 20354.44
$ws.Range("C3").Value = 150139.79
$ws.Range("D3").Value = 187082.86
$ws.Range("E3").Value = 343644.5
$ws.Range("F3").Value = 453755.79

$ws.Range("B4").Value = 65838.63
$ws.Range("C4").Value = 164999.57
$ws.Range("D4").Value = 193765.12
$ws.Range("E4").Value = 323694.55
$ws.Range("F4").Value = 158167.64

$ws.Range("B5").Value = 46161.32
$ws.Range("C5").Value = 130298.05
$ws.Range("D5").Value = 215432.77
$ws.Range("E5").Value = 396044.77

$ws.Range("B6").Value = 47815.15
$ws.Range("C6").Value = 115915.98
$ws.Range("D6").Value = 220684.43
$ws.Range("E6").Value = 378096.36

$ws.Range("B7").Value = 64428.92
$ws.Range("C7").Value = 141564.4
$ws.Range("D7").Value = 246532.18
$ws.Range("E7").Value = 403206.55

$ws.Range("B8").Value = 72571.03
$ws.Range("C8").Value = 136354.87
$ws.Range("D8").Value = 223998.98
$ws.Range("E8").Value = 389318.8

$ws.Range("B9").Value = 116043.68
$ws.Range("C9").Value = 144340.4
$ws.Range("D9").Value = 267309.71
$ws.Range("E9").Value = 381723.02

$ws.Range("B10").Value = 112241.8
$ws.Range("C10").Value = 153181.38
$ws.Range("D10").Value = 292606.09
$ws.Range("E10").Value = 386980.51

$ws.Range("B11").Value = 117919.45
$ws.Range("C11").Value = 134868.15
$ws.Range("D11").Value = 331822.65
$ws.Range("E11").Value = 418014.5

$ws.Range("B12").Value = 129704.99
$ws.Range("C12").Value = 151292.14
$ws.Range("D12").Value = 249653.56
$ws.Range("E12").Value = 309538.28

$ws.Range("B13").Value = 106870.37
$ws.Range("C13").Value = 201369.38
$ws.Range("D13").Value = 251128.74
$ws.Range("E13").Value = 381502.92
